$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in / correct the "Absent" column (H) so that it reflects the
# "Total Attendance Count" column (D): Absent = 1 when no attendance was
# recorded that day (D = 0), Absent = 0 when attendance was recorded (D = 1).
# This forms the consolidated report for the remaining blank/incorrect rows.

$ws.Range("H3").Value = 1
$ws.Range("H5").Value = 0
$ws.Range("H9").Value = 1
$ws.Range("H12").Value = 0
$ws.Range("H13").Value = 1
$ws.Range("H14").Value = 0
$ws.Range("H19").Value = 1
$ws.Range("H21").Value = 0
